$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the new row 20, cloning the constant columns from row 19 ---
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 44376
$ws.Range("D20").NumberFormat = $ws.Range("D19").NumberFormat
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = 100112052
$ws.Range("G20").Value = "Albahaca"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 270
$ws.Range("K20").Value = 2400
$ws.Range("L20").Value = 2500
$ws.Range("M20").Value = 2437
$ws.Range("N20").Value = "$/paquete"
$ws.Range("O20").Value = "Región de Arica y Parinacota"
$ws.Range("P20").Value = 2437
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Hortaliza"

# --- Update weekly Fecha / Volumen / Precio columns for existing rows (shuffled weekly snapshot) ---
$ws.Range("D2").Value = 44326
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 2700
$ws.Range("L2").Value = 2800
$ws.Range("M2").Value = 2750
$ws.Range("P2").Value = 2750

$ws.Range("D3").Value = 44349
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 2800
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 2900
$ws.Range("P3").Value = 2900

$ws.Range("D4").Value = 44280
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 1400
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = 1450
$ws.Range("P4").Value = 1450

$ws.Range("D5").Value = 44260
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 900
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = 950
$ws.Range("P5").Value = 950

$ws.Range("D6").Value = 44330
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 2800
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = 2900
$ws.Range("P6").Value = 2900

$ws.Range("D7").Value = 44250
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 1200
$ws.Range("M7").Value = 1100
$ws.Range("P7").Value = 1100

$ws.Range("D8").Value = 44292
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 1800
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 1900
$ws.Range("P8").Value = 1900

$ws.Range("D9").Value = 44224
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 750
$ws.Range("L9").Value = 800
$ws.Range("M9").Value = 775
$ws.Range("P9").Value = 775

$ws.Range("D10").Value = 44362
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 2800
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 2900
$ws.Range("P10").Value = 2900

$ws.Range("D11").Value = 44249
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 900
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = 950
$ws.Range("P11").Value = 950

$ws.Range("D12").Value = 44365
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 2400
$ws.Range("L12").Value = 2500
$ws.Range("M12").Value = 2450
$ws.Range("P12").Value = 2450

$ws.Range("D13").Value = 44274
$ws.Range("J13").Value = 250
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 1200
$ws.Range("M13").Value = 1100
$ws.Range("P13").Value = 1100

$ws.Range("D14").Value = 44417
$ws.Range("J14").Value = 250
$ws.Range("K14").Value = 4000
$ws.Range("L14").Value = 4500
$ws.Range("M14").Value = 4250
$ws.Range("P14").Value = 4250

$ws.Range("D15").Value = 44432
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 2300
$ws.Range("L15").Value = 2500
$ws.Range("M15").Value = 2400
$ws.Range("P15").Value = 2400

$ws.Range("D16").Value = 44302
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 950
$ws.Range("P16").Value = 950

$ws.Range("D17").Value = 44313
$ws.Range("J17").Value = 250
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 1000
$ws.Range("M17").Value = 950
$ws.Range("P17").Value = 950

$ws.Range("D18").Value = 44442
$ws.Range("J18").Value = 240
$ws.Range("K18").Value = 2300
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = 2400
$ws.Range("P18").Value = 2400

$ws.Range("D19").Value = 44435
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 2300
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = 2400
$ws.Range("P19").Value = 2400

